$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InvalidLogin")

# --- Insert a new test case row at row 5: "standard_user@" / "secret_sauce" ---
# (an extra '@' typo in the username -> still just a username/password mismatch)
$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = "standard_user@"
$ws.Range("B5").Value = "secret_sauce"
$ws.Range("C5").Value = "Epic sadface: Username and password do not match any user in this service"

# --- Append two more new test cases at the bottom of the table (rows 12 & 13) ---
# (the earlier insert at row 5 pushed the old last data row from 10 down to 11)

# Row 12: correct username, "secret_sauce#" (typo'd password) -> mismatch error
$ws.Rows.Item(11).Copy()
$ws.Rows.Item(12).Insert()
$ws.Range("A12").Value = "standard_user"
$ws.Range("B12").Value = "secret_sauce#"
$ws.Range("C12").Value = "Epic sadface: Username and password do not match any user in this service"

# Row 13: "locked_out_user" / correct password -> locked-out error
$ws.Rows.Item(12).Copy()
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Value = "locked_out_user"
$ws.Range("B13").Value = "secret_sauce"
$ws.Range("C13").Value = "Epic sadface: Sorry, this user has been locked out."
$ws.Range("C13").Font.Size = 10
